{"js": "// Update the multiplication-problem table: each \"A\u00d7B=\" cell gets replaced\n// with a new \"C\u00d7D=\" expression, one-for-one, per the commit's diff.\nconst replacements = [\n  [\"40\u00d776=\", \"90\u00d793=\"],\n  [\"65\u00d785=\", \"17\u00d741=\"],\n  [\"47\u00d792=\", \"61\u00d734=\"],\n  [\"46\u00d774=\", \"99\u00d743=\"],\n  [\"78\u00d760=\", \"91\u00d768=\"],\n  [\"52\u00d774=\", \"39\u00d757=\"],\n  [\"40\u00d779=\", \"89\u00d731=\"],\n  [\"77\u00d788=\", \"68\u00d773=\"],\n  [\"33\u00d745=\", \"29\u00d777=\"],\n  [\"36\u00d757=\", \"92\u00d754=\"],\n  [\"61\u00d766=\", \"18\u00d733=\"],\n  [\"27\u00d777=\", \"16\u00d795=\"],\n  [\"70\u00d736=\", \"44\u00d747=\"],\n  [\"22\u00d750=\", \"72\u00d740=\"],\n  [\"96\u00d790=\", \"87\u00d725=\"],\n  [\"81\u00d753=\", \"21\u00d727=\"],\n  [\"56\u00d764=\", \"46\u00d712=\"],\n  [\"85\u00d728=\", \"69\u00d769=\"],\n  [\"40\u00d711=\", \"88\u00d779=\"],\n  [\"95\u00d732=\", \"30\u00d780=\"],\n  [\"12\u00d765=\", \"97\u00d766=\"],\n  [\"40\u00d731=\", \"60\u00d788=\"],\n  [\"61\u00d742=\", \"98\u00d795=\"],\n  [\"86\u00d741=\", \"37\u00d766=\"],\n  [\"54\u00d799=\", \"57\u00d739=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the multiplication-problem table: each \"A\u00d7B=\" cell gets replaced\n# with a new \"C\u00d7D=\" expression, one-for-one, per the commit's diff.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"40\u00d776=\"; New = \"90\u00d793=\" },\n    @{ Old = \"65\u00d785=\"; New = \"17\u00d741=\" },\n    @{ Old = \"47\u00d792=\"; New = \"61\u00d734=\" },\n    @{ Old = \"46\u00d774=\"; New = \"99\u00d743=\" },\n    @{ Old = \"78\u00d760=\"; New = \"91\u00d768=\" },\n    @{ Old = \"52\u00d774=\"; New = \"39\u00d757=\" },\n    @{ Old = \"40\u00d779=\"; New = \"89\u00d731=\" },\n    @{ Old = \"77\u00d788=\"; New = \"68\u00d773=\" },\n    @{ Old = \"33\u00d745=\"; New = \"29\u00d777=\" },\n    @{ Old = \"36\u00d757=\"; New = \"92\u00d754=\" },\n    @{ Old = \"61\u00d766=\"; New = \"18\u00d733=\" },\n    @{ Old = \"27\u00d777=\"; New = \"16\u00d795=\" },\n    @{ Old = \"70\u00d736=\"; New = \"44\u00d747=\" },\n    @{ Old = \"22\u00d750=\"; New = \"72\u00d740=\" },\n    @{ Old = \"96\u00d790=\"; New = \"87\u00d725=\" },\n    @{ Old = \"81\u00d753=\"; New = \"21\u00d727=\" },\n    @{ Old = \"56\u00d764=\"; New = \"46\u00d712=\" },\n    @{ Old = \"85\u00d728=\"; New = \"69\u00d769=\" },\n    @{ Old = \"40\u00d711=\"; New = \"88\u00d779=\" },\n    @{ Old = \"95\u00d732=\"; New = \"30\u00d780=\" },\n    @{ Old = \"12\u00d765=\"; New = \"97\u00d766=\" },\n    @{ Old = \"40\u00d731=\"; New = \"60\u00d788=\" },\n    @{ Old = \"61\u00d742=\"; New = \"98\u00d795=\" },\n    @{ Old = \"86\u00d741=\"; New = \"37\u00d766=\" },\n    @{ Old = \"54\u00d799=\"; New = \"57\u00d739=\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $find.Execute($pair.Old, $false, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2) | Out-Null\n}\n\nWrite-Output \"done\"\n"}
